# Generate Report for Handback
# Refresh the timestamp strings recorded on the handback-status report.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date for the first file (shown on both the
# Overview rollup and the de-de detail sheet's Correspond Handoff Datetime).
$overview.Range("G2").Value = "2016-09-04 09:09:23"
$dede.Range("H2").Value = "2016-09-04 09:09:23"

# zh-cn detail sheet: handoff + handback datetimes for the first file.
$zhcn.Range("H2").Value = "2016-09-04 09:09:19"
$zhcn.Range("K2").Value = "2016-09-04 09:09:37"

# de-de detail sheet: handback datetime for the first file.
$dede.Range("K2").Value = "2016-09-04 09:09:44"
